# Loan and Deposit Balance Report
# Include Credit officer on Loan Accounts and Deposit Accounts reports
#
# Insert a new "LOAN OFFICER" column right after OFFICE (before CLIENT ID),
# shifting CLIENT ID / NAME / TYPE / ACCRUED INTEREST / BALANCE / STATUS
# one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column C (CLIENT ID).
$ws.Columns("C").Insert()

# New header cell + text.
$ws.Range("C2").Value = "LOAN OFFICER"

# Match the header formatting used by the other header cells (bold, border, centered).
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Match the column width used by its neighbouring header column.
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth

# Reflect where the cursor ended up after the edit.
$ws.Range("D11").Select()
